$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = 'Data Analyst'
$ws.Range('B2').Value = 'Confidential'
$ws.Range('D2').Value = 'https://wuzzuf.net/jobs/p/ewn4cDHL1bUW-Data-Analyst-Giza-Egypt'
$ws.Range('A3').Value = 'Centro Data Program'
$ws.Range('B3').Value = 'Centro'
$ws.Range('C3').Value = 'Maadi, Cairo, Egypt'
$ws.Range('D3').Value = 'https://wuzzuf.net/internship/j86RYZ2Qg0Wn-Centro-Data-Program-Centro-Cairo-Egypt'
$ws.Range('E3').Value = 'Internship'
$ws.Range('A4').Value = 'Data Analyst'
$ws.Range('B4').Value = 'CRIF EGYPT'
$ws.Range('C4').Value = 'Heliopolis, Cairo, Egypt'
$ws.Range('D4').Value = 'https://wuzzuf.net/jobs/p/AXcIDiTsKSWn-Data-Analyst-CRIF-EGYPT-Cairo-Egypt'
$ws.Range('A5').Value = 'Data Analyst Team Leader'
$ws.Range('B5').Value = '3Sixty'
$ws.Range('C5').Value = 'Mohandessin, Giza, Egypt'
$ws.Range('D5').Value = 'https://wuzzuf.net/jobs/p/HN6h2CLSlb04-Data-Analyst-Team-Leader-3Sixty-Giza-Egypt'
$ws.Range('E5').Value = 'Full Time'
$ws.Range('F5').Value = 'Hybrid'
$ws.Range('A6').Value = 'Business Data Analyst'
$ws.Range('B6').Value = 'Fawry Plus'
$ws.Range('C6').Value = 'Smart Village, Giza, Egypt'
$ws.Range('D6').Value = 'https://wuzzuf.net/jobs/p/kvMsx40qvLZl-Business-Data-Analyst-Fawry-Plus-Giza-Egypt'
$ws.Range('A7').Value = 'Big Data Consultant'
$ws.Range('B7').Value = 'Exco Egypt'
$ws.Range('C7').Value = 'Nasr City, Cairo, Egypt'
$ws.Range('D7').Value = 'https://wuzzuf.net/jobs/p/egjM6VMefvKW-Big-Data-Consultant-Exco-Egypt-Cairo-Egypt'
$ws.Range('A8').Value = 'Data/Business Analyst Intern'
$ws.Range('B8').Value = 'kafiil'
$ws.Range('D8').Value = 'https://wuzzuf.net/internship/NXl1N6zwtmSd-DataBusiness-Analyst-Intern-kafiil-Cairo-Egypt'
$ws.Range('E8').Value = 'Internship'
$ws.Range('A9').Value = 'Financial Data Analyst'
$ws.Range('B9').Value = 'Influence Communication'
$ws.Range('C9').Value = 'Maadi, Cairo, Egypt'
$ws.Range('D9').Value = 'https://wuzzuf.net/jobs/p/pDJ6Rtkhmt6y-Financial-Data-Analyst-Influence-Communication-Cairo-Egypt'
$ws.Range('A10').Value = 'Data Management Lead'
$ws.Range('B10').Value = 'BBI-Consultancy'
$ws.Range('C10').Value = 'Cairo, Egypt'
$ws.Range('D10').Value = 'https://wuzzuf.net/jobs/p/KhhiGntCocD4-Data-Management-Lead-BBI-Consultancy-Cairo-Egypt'
$ws.Range('A11').Value = 'Data Analytics Engineer'
$ws.Range('B11').Value = 'Erada'
$ws.Range('C11').Value = 'Maadi, Cairo, Egypt'
$ws.Range('D11').Value = 'https://wuzzuf.net/jobs/p/M0QOPoyJQ6Ia-Data-Analytics-Engineer-Erada-Cairo-Egypt'
$ws.Range('B12').Value = 'Talaat Moustafa Group'
$ws.Range('C12').Value = 'Dokki, Giza, Egypt'
$ws.Range('D12').Value = 'https://wuzzuf.net/jobs/p/n8eZCC5B0pgZ-Data-Analyst-Talaat-Moustafa-Group-Giza-Egypt'
$ws.Range('A13').Value = 'BI & Data Analytics Consultant'
$ws.Range('B13').Value = 'Mantrac'
$ws.Range('C13').Value = 'Alexandria, Egypt'
$ws.Range('D13').Value = 'https://wuzzuf.net/jobs/p/B2lwaYBkrVx3-BI-Data-Analytics-Consultant-Mantrac-Alexandria-Egypt'
$ws.Range('B14').Value = 'Mazaya'
$ws.Range('C14').Value = 'New Cairo, Cairo, Egypt'
$ws.Range('D14').Value = 'https://wuzzuf.net/jobs/p/HaCk4dcFubh3-Data-Analyst-Mazaya-Cairo-Egypt'
$ws.Range('A15').Value = 'Data Analyst'
$ws.Range('B15').Value = 'Rotem SRS'
$ws.Range('C15').Value = 'Cairo, Egypt'
$ws.Range('D15').Value = 'https://wuzzuf.net/jobs/p/dZtsSzWqvAqg-Data-Analyst-Rotem-SRS-Cairo-Egypt'
$ws.Range('B16').Value = 'Al Watania Poultry'
$ws.Range('C16').Value = 'Smart Village, Giza, Egypt'
$ws.Range('D16').Value = 'https://wuzzuf.net/jobs/p/HzJ6GhNp7osP-Data-Analyst-Al-Watania-Poultry-Giza-Egypt'
